$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "date" column (F) values forward by 3 days for rows 2-7
$ws.Range("F2").Value = 44588
$ws.Range("F3").Value = 44587
$ws.Range("F4").Value = 44586
$ws.Range("F5").Value = 44585
$ws.Range("F6").Value = 44584
$ws.Range("F7").Value = 44583
